$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "3"
$ws.Range("G2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "3"
$ws.Range("G3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.427"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "3"
$ws.Range("G4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05987"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "3"
$ws.Range("G5").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "3"
$ws.Range("G6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.516"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "3"
$ws.Range("G7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8152"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "3"
$ws.Range("G8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9232"
$ws.Range("D9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "3"
$ws.Range("G9").Style = "Normal"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1413"
$ws.Range("D10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "3"
$ws.Range("G10").Style = "Normal"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07374"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3"
$ws.Range("G11").Style = "Normal"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03230"
$ws.Range("D12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "3"
$ws.Range("G12").Style = "Normal"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03056"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "3"
$ws.Range("G13").Style = "Normal"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09355"
$ws.Range("D14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "3"
$ws.Range("G14").Style = "Normal"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.845"
$ws.Range("D15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "3"
$ws.Range("G15").Style = "Normal"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001559"
$ws.Range("D16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "3"
$ws.Range("G16").Style = "Normal"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04682"
$ws.Range("D17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "3"
$ws.Range("G17").Style = "Normal"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005939"
$ws.Range("D18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "3"
$ws.Range("G18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006120"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "3"
$ws.Range("G19").Style = "Normal"

$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005017"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "3"
$ws.Range("G20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009839"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "3"
$ws.Range("G21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00007795"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "3"
$ws.Range("G22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.624"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "3"
$ws.Range("G23").Style = "Normal"

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "3"
$ws.Range("G24").Style = "Normal"

$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "3"
$ws.Range("G25").Style = "Normal"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "3"
$ws.Range("G26").Style = "Normal"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "3"
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "3"
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "3"
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "3"
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "3"
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "3"
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "3"
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "3"
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "3"
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "3"
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "3"
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "3"
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "3"
$ws.Range("G39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03924"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "3"
$ws.Range("G40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006198"
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "3"
$ws.Range("G41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1077"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "3"
$ws.Range("G42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002618"
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "3"
$ws.Range("G43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007105"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "3"
$ws.Range("G44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005242"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "3"
$ws.Range("G45").Style = "Normal"

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "3"
$ws.Range("G46").Style = "Normal"

$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "3"
$ws.Range("G47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9099"
$ws.Range("D48").Style = "Normal"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "3"
$ws.Range("G48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002298"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "3"
$ws.Range("G49").Style = "Normal"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "3"
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "3"
$ws.Range("G51").Style = "Normal"

